# Auto-generated: applies scheduled-runner market-data refresh to Sheets/Bahamut_Profits.xlsx
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and derived Leve Profit columns
# for the affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 2654.7896
$ws.Range("I132").Value = 2725.9412
$ws.Range("K132").Value = 8177.823600000001
$ws.Range("M132").Value = -5647.823600000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 6892.7793
$ws.Range("I32").Value = 5098.403
$ws.Range("K32").Value = 5098.403
$ws.Range("M32").Value = -4811.403

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 973
$ws.Range("I74").Value = 1215.2941
$ws.Range("J74").Value = 561.1
$ws.Range("K74").Value = 1215.2941
$ws.Range("L74").Value = 561.1
$ws.Range("M74").Value = -341.2941000000001
$ws.Range("N74").Value = -2309.1

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 973
$ws.Range("I77").Value = 1215.2941
$ws.Range("J77").Value = 561.1
$ws.Range("K77").Value = 6076.4705
$ws.Range("L77").Value = 2805.5
$ws.Range("M77").Value = -1708.4705
$ws.Range("N77").Value = -11541.5

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 726.1111
$ws.Range("I97").Value = 579
$ws.Range("J97").Value = 1020.3333
$ws.Range("K97").Value = 579
$ws.Range("L97").Value = 1020.3333
$ws.Range("M97").Value = -83
$ws.Range("N97").Value = -2012.3333

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 3334
$ws.Range("I102").Value = 3369.889
$ws.Range("J102").Value = 3011
$ws.Range("K102").Value = 3369.889
$ws.Range("L102").Value = 3011
$ws.Range("M102").Value = -1747.889
$ws.Range("N102").Value = -6255

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1108.3334
$ws.Range("I110").Value = 950
$ws.Range("J110").Value = 1187.5
$ws.Range("K110").Value = 950
$ws.Range("L110").Value = 1187.5
$ws.Range("M110").Value = 1095
$ws.Range("N110").Value = -5277.5

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 1725.0513
$ws.Range("I132").Value = 1299.9678
$ws.Range("J132").Value = 3372.25
$ws.Range("K132").Value = 3899.9034
$ws.Range("L132").Value = 10116.75
$ws.Range("M132").Value = -1369.9034
$ws.Range("N132").Value = -15176.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 36564.5
$ws.Range("I20").Value = 87013
$ws.Range("J20").Value = 9399.923000000001
$ws.Range("K20").Value = 87013
$ws.Range("L20").Value = 9399.923000000001
$ws.Range("M20").Value = -86766
$ws.Range("N20").Value = -9893.923000000001

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 52634464
$ws.Range("I86").Value = 62502176
$ws.Range("J86").Value = 6666.6665
$ws.Range("K86").Value = 62502176
$ws.Range("L86").Value = 6666.6665
$ws.Range("M86").Value = -62501053
$ws.Range("N86").Value = -8912.666499999999

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 52634464
$ws.Range("I89").Value = 62502176
$ws.Range("J89").Value = 6666.6665
$ws.Range("K89").Value = 312510880
$ws.Range("L89").Value = 33333.3325
$ws.Range("M89").Value = -312505264
$ws.Range("N89").Value = -44565.3325

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 878.78125
$ws.Range("I94").Value = 847.0333000000001
$ws.Range("K94").Value = 847.0333000000001
$ws.Range("M94").Value = -396.0333000000001

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 1912.3
$ws.Range("I105").Value = 2051.6155
$ws.Range("J105").Value = 1006.75
$ws.Range("K105").Value = 2051.6155
$ws.Range("L105").Value = 1006.75
$ws.Range("M105").Value = -304.6154999999999
$ws.Range("N105").Value = -4500.75

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 2301.0715
$ws.Range("I107").Value = 2021.409
$ws.Range("J107").Value = 3326.5
$ws.Range("K107").Value = 2021.409
$ws.Range("L107").Value = 3326.5
$ws.Range("M107").Value = -101.4090000000001
$ws.Range("N107").Value = -7166.5

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2895.2122
$ws.Range("I134").Value = 2258.12
$ws.Range("J134").Value = 4886.125
$ws.Range("K134").Value = 6774.36
$ws.Range("L134").Value = 14658.375
$ws.Range("M134").Value = -4239.36
$ws.Range("N134").Value = -19728.375

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 1800
$ws.Range("I16").Value = 1800
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1513

# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 1079.9333
$ws.Range("I22").Value = 1300.091
$ws.Range("J22").Value = 474.5
$ws.Range("K22").Value = 1300.091
$ws.Range("L22").Value = 474.5
$ws.Range("M22").Value = -950.0909999999999
$ws.Range("N22").Value = -1174.5

# Row 35: Storm of Swords | Elm Macuahuitl
$ws.Range("H35").Value = 370
$ws.Range("I35").Value = 370
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 370
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -76

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2403
$ws.Range("I58").Value = 2403
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2403
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -2200

# Row 94: Beech, Please | Beech Lumber
$ws.Range("H94").Value = 4170.143
$ws.Range("I94").Value = 3627.5715
$ws.Range("J94").Value = 4351
$ws.Range("K94").Value = 3627.5715
$ws.Range("L94").Value = 4351
$ws.Range("M94").Value = -3176.5715
$ws.Range("N94").Value = -5253

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 3672.818
$ws.Range("I105").Value = 4578
$ws.Range("J105").Value = 2918.5
$ws.Range("K105").Value = 4578
$ws.Range("L105").Value = 2918.5
$ws.Range("M105").Value = -2831
$ws.Range("N105").Value = -6412.5

# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 326.54166
$ws.Range("I107").Value = 226.58333
$ws.Range("J107").Value = 426.5
$ws.Range("K107").Value = 226.58333
$ws.Range("L107").Value = 426.5
$ws.Range("M107").Value = 1693.41667
$ws.Range("N107").Value = -4266.5

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 1800
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 370

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2408.5454
$ws.Range("I132").Value = 1170.1538
$ws.Range("J132").Value = 4197.3335
$ws.Range("K132").Value = 3510.4614
$ws.Range("L132").Value = 12592.0005
$ws.Range("M132").Value = -980.4614000000001
$ws.Range("N132").Value = -17652.0005

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2403
$ws.Range("I136").Value = 2403
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7209
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -4659

$ws = $wb.Worksheets.Item("CUL")
# Row 103: West Meats East | Nomad Meat Pie
$ws.Range("H103").Value = 1172.8
$ws.Range("I103").Value = 216
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 648
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 231
$ws.Range("N103").Value = -16758

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 631.8
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 574.7917
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 1724.3751
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -6064.3751

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 60.875
$ws.Range("I2").Value = 47.083332
$ws.Range("J2").Value = 102.25
$ws.Range("K2").Value = 47.083332
$ws.Range("L2").Value = 102.25
$ws.Range("M2").Value = 65.916668
$ws.Range("N2").Value = -328.25

# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 4165.3
$ws.Range("I70").Value = 4139.56
$ws.Range("K70").Value = 4139.56
$ws.Range("M70").Value = -3869.56

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 4165.3
$ws.Range("I73").Value = 4139.56
$ws.Range("K73").Value = 4139.56
$ws.Range("M73").Value = -3203.56

# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 4065.1177
$ws.Range("I80").Value = 4300
$ws.Range("J80").Value = 3800.875
$ws.Range("K80").Value = 4300
$ws.Range("L80").Value = 3800.875
$ws.Range("M80").Value = -3302
$ws.Range("N80").Value = -5796.875

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 4065.1177
$ws.Range("I83").Value = 4300
$ws.Range("J83").Value = 3800.875
$ws.Range("K83").Value = 21500
$ws.Range("L83").Value = 19004.375
$ws.Range("M83").Value = -16508
$ws.Range("N83").Value = -28988.375

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 1022.65515
$ws.Range("I97").Value = 1086.72
$ws.Range("J97").Value = 622.25
$ws.Range("K97").Value = 1086.72
$ws.Range("L97").Value = 622.25
$ws.Range("M97").Value = -590.72
$ws.Range("N97").Value = -1614.25

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 1816.8
$ws.Range("I102").Value = 1816.8
$ws.Range("K102").Value = 1816.8
$ws.Range("M102").Value = -194.8

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 1974.8182
$ws.Range("I113").Value = 662.5
$ws.Range("J113").Value = 2724.7144
$ws.Range("K113").Value = 662.5
$ws.Range("L113").Value = 2724.7144
$ws.Range("M113").Value = 1507.5
$ws.Range("N113").Value = -7064.7144

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2982.825
$ws.Range("I132").Value = 2686.3
$ws.Range("J132").Value = 3872.4
$ws.Range("K132").Value = 8058.900000000001
$ws.Range("L132").Value = 11617.2
$ws.Range("M132").Value = -5528.900000000001
$ws.Range("N132").Value = -16677.2

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 1540.375
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1540.375
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 1540.375
$ws.Range("N22").Value = -2130.375

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 1540.375
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1540.375
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 1540.375
$ws.Range("N27").Value = -1754.375

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Range("H55").Value = 191.16667
$ws.Range("I55").Value = 199.4
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 199.4
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = -26.40000000000001
$ws.Range("N55").Value = -496

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 1214.7142
$ws.Range("I93").Value = 1139.2727
$ws.Range("J93").Value = 1491.3334
$ws.Range("K93").Value = 1139.2727
$ws.Range("L93").Value = 1491.3334
$ws.Range("M93").Value = 108.7273
$ws.Range("N93").Value = -3987.3334

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 6612.826
$ws.Range("I122").Value = 9862.154
$ws.Range("K122").Value = 29586.462
$ws.Range("M122").Value = -27136.462

$ws = $wb.Worksheets.Item("WVR")
# Row 51: After the Smock-down | Linen Smock
$ws.Range("H51").Value = 40000
$ws.Range("J51").Value = 40000
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41020

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 368.83334
$ws.Range("I126").Value = 426.8421
$ws.Range("J126").Value = 148.4
$ws.Range("K126").Value = 1280.5263
$ws.Range("L126").Value = 445.2
$ws.Range("M126").Value = 1189.4737
$ws.Range("N126").Value = -5385.2

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2466
$ws.Range("I132").Value = 2115.3635
$ws.Range("J132").Value = 3108.8333
$ws.Range("K132").Value = 6346.0905
$ws.Range("L132").Value = 9326.499899999999
$ws.Range("M132").Value = -3816.0905
$ws.Range("N132").Value = -14386.4999

